$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E ("calendly link") for the existing rows -----------------
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "calendly link"
$ws.Range("E2").Value = "https://calendly.com/nick-griffiths-22/strategy-meeting-clone"
$ws.Range("E3").Value = "https://calendly.com/nick-griffiths-22/strategy-meeting-clone"
$ws.Range("E1").ColumnWidth = 53

# --- New rows 4-6, cloned from row 3 (Yoda / yoda.png / study hard you must)
$names = @("Yoda 2", "Yoda 3", "Yoda 4")
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 4 + $i

    # Copy row 3 values + formats into the new row, then rename column A.
    $ws.Range("A3:E3").Copy()
    $ws.Range("A$row").PasteSpecial(-4163)
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $names[$i]

    # Give the new row's email cell its own hyperlink (mirrors D3).
    $ws.Hyperlinks.Add($ws.Range("D$row"), "mailto:yoda@email.com") | Out-Null

    # Hyperlinks.Add() re-stamps the cell style; restore the shared
    # "Hyperlink + wrap text" formatting used by D2/D3.
    $ws.Range("D3").Copy()
    $ws.Range("D$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# --- View bits (matches the committed sheet state) -------------------------
$ws.Range("A6").Select() | Out-Null
